# Adds a new emission-flow row for "Carbon dioxide, non-fossil, resource
# correction" right after the existing "Carbon dioxide, in air" row, pushing
# every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3 (shifts rows 3..218 down to 4..219).
$ws.Rows.Item(3).Insert()

# Populate the new row's data.
$ws.Cells.Item(3, 1).Value() = "Carbon dioxide, non-fossil, resource correction"
$ws.Cells.Item(3, 2).Value() = "natural resource::in air"
$ws.Cells.Item(3, 3).Value() = 1

# Reset the view: scroll back to the top and move the selection to E13.
$ws.Activate()
$ws.Range("E13").Select()
